$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 235, shifting existing rows 235-341 down to 236-342
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with the new data record
$ws.Range("A235").Value = 11
$ws.Range("B235").Value = "Vega Monumental Concepción"
$ws.Range("C235").Value = "Bíobío"
$ws.Range("D235").Value = 44813
$ws.Range("E235").Value = 8
$ws.Range("F235").Value = "Fruta"
$ws.Range("G235").Value = 100102
$ws.Range("H235").Value = "Cítricos"
$ws.Range("I235").Value = 100102005
$ws.Range("J235").Value = "Naranja"
$ws.Range("K235").Value = "Lane Late"
$ws.Range("L235").Value = "Primera"
$ws.Range("M235").Value = 220
$ws.Range("N235").Value = 6000
$ws.Range("O235").Value = 6500
$ws.Range("P235").Value = 6273
$ws.Range("Q235").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R235").Value = "Región de O'Higgins"
$ws.Range("S235").Value = 418
$ws.Range("T235").Value = 15
